# This script applies targeted numeric corrections to the stock report
# worksheet, matching the authoritative re-count of quantities / values
# for the affected SKU rows (and the two duplicate-item row pairs whose
# rows were swapped), and propagates the resulting differences up into
# each "Sub Total:" / "Grand Total:" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("F21").Value = 142
    $ws.Range("G21").Value = 3645.14
    $ws.Range("B32").Value = 12174.01
    $ws.Range("F64").Value = 108
    $ws.Range("G64").Value = 8766.360000000001
    $ws.Range("F71").Value = 291
    $ws.Range("G71").Value = 18536.7
    $ws.Range("F73").Value = 65
    $ws.Range("G73").Value = 5133.7
    $ws.Range("F85").Value = 133
    $ws.Range("G85").Value = 17924.41
    $ws.Range("B90").Value = 161040.91
    $ws.Range("F144").Value = 929
    $ws.Range("G144").Value = 7850.05
    $ws.Range("B147").Value = 11969.03
    $ws.Range("F150").Value = 22
    $ws.Range("G150").Value = 1022.78
    $ws.Range("B156").Value = 27989.81
    $ws.Range("F182").Value = 17
    $ws.Range("G182").Value = 3041.47
    $ws.Range("B192").Value = 64973
    $ws.Range("E192").Value = 35.4
    $ws.Range("F192").Value = 0
    $ws.Range("G192").Value = 0
    $ws.Range("B193").Value = 48706
    $ws.Range("E193").Value = 39.8
    $ws.Range("F193").Value = -144
    $ws.Range("G193").Value = -4795.2
    $ws.Range("F205").Value = 15
    $ws.Range("G205").Value = 5657.1
    $ws.Range("B216").Value = 29710.66
    $ws.Range("F249").Value = 129
    $ws.Range("G249").Value = 17778.78
    $ws.Range("B260").Value = 161875.3
    $ws.Range("F280").Value = 126
    $ws.Range("G280").Value = 21311.64
    $ws.Range("F302").Value = 24
    $ws.Range("G302").Value = 5061.36
    $ws.Range("B304").Value = 158556.45
    $ws.Range("F320").Value = 33
    $ws.Range("G320").Value = 2265.45
    $ws.Range("F321").Value = 42
    $ws.Range("G321").Value = 2306.64
    $ws.Range("B330").Value = 24775.55
    $ws.Range("F345").Value = 26
    $ws.Range("G345").Value = 1596.66
    $ws.Range("B346").Value = 22897.77
    $ws.Range("B382").Value = 53263
    $ws.Range("E382").Value = 15.29
    $ws.Range("F382").Value = -309
    $ws.Range("G382").Value = -3958.29
    $ws.Range("B383").Value = 65066
    $ws.Range("E383").Value = 13.61
    $ws.Range("F383").Value = 90
    $ws.Range("G383").Value = 1152.9
    $ws.Range("B391").Value = 64927
    $ws.Range("E391").Value = 17.26
    $ws.Range("F391").Value = 106
    $ws.Range("G391").Value = 1719.32
    $ws.Range("B392").Value = 45718
    $ws.Range("E392").Value = 19.38
    $ws.Range("F392").Value = -294
    $ws.Range("G392").Value = -4768.68
    $ws.Range("B396").Value = 45709
    $ws.Range("E396").Value = 15.69
    $ws.Range("F396").Value = -300
    $ws.Range("G396").Value = -3945
    $ws.Range("B397").Value = 64925
    $ws.Range("E397").Value = 13.97
    $ws.Range("F397").Value = 111
    $ws.Range("G397").Value = 1459.65
    $ws.Range("F472").Value = 42
    $ws.Range("G472").Value = 4643.94
    $ws.Range("B476").Value = 11964.4
    $ws.Range("F526").Value = 50
    $ws.Range("G526").Value = 5197
    $ws.Range("B532").Value = 132392.66
    $ws.Range("F599").Value = 34
    $ws.Range("G599").Value = 1461.66
    $ws.Range("F604").Value = 16
    $ws.Range("G604").Value = 911.84
    $ws.Range("B605").Value = 11669.88
    $ws.Range("F621").Value = 1213
    $ws.Range("G621").Value = 197852.43
    $ws.Range("B628").Value = 341527.5
    $ws.Range("F641").Value = 124
    $ws.Range("G641").Value = 19735.84
    $ws.Range("B646").Value = 40871.38
    $ws.Range("B647").Value = 2265633.81
    $ws.Range("B648").Value = 2265633.81
